$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.277.23'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.458.30'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.23%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.455.91'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.90%  '
$ws.Range('E11').Value = '  +4.65%  '
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.055.67'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.10'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.54%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000195'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.300.03'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.450.88'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.44'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '387.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.550'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.24%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000124'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +20.07%  '
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +10.40%  '
$ws.Range('E31').Value = '  +9.81%  '
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.60'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.77'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.14'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.51'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.76'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0779'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.25%  '
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.51'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.923.02'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0323'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.51'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.94%  '
$ws.Range('E45').Value = '  +3.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.770'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +8.80%  '
$ws.Range('E48').Value = '  +2.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.20'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +15.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.108'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.87%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.868'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.83%  '
